# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
# Rows 128/129 swap their match records (cols B, F:AC).
# Rows 148/149/150 rotate their match records: 148<-150(old), 149<-148(old), 150<-149(old).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128
$ws.Range("B128").Value = 7462738
$ws.Range("F128").Value = 'Vaca Diez'
$ws.Range("G128").Value = 'The Strongest'
$ws.Range("H128").Value = 2
$ws.Range("I128").Value = 2
$ws.Range("J128").Value = 'D'
$ws.Range("K128").Value = 4
$ws.Range("L128").Value = 4
$ws.Range("M128").Value = 1.666
$ws.Range("N128").Value = 4
$ws.Range("O128").Value = 3.8
$ws.Range("P128").Value = 1.75
$ws.Range("Q128").Value = 0.75
$ws.Range("R128").Value = 1.8
$ws.Range("S128").Value = 2
$ws.Range("T128").Value = 3
$ws.Range("U128").Value = 1.925
$ws.Range("V128").Value = 1.875
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = 2.8
$ws.Range("Y128").Value = -1
$ws.Range("Z128").Value = 0.8
$ws.Range("AA128").Value = -1
$ws.Range("AB128").Value = 0.925
$ws.Range("AC128").Value = -1

# Row 129
$ws.Range("B129").Value = 7462542
$ws.Range("F129").Value = 'Always Ready'
$ws.Range("G129").Value = 'Royal Pari FC'
$ws.Range("H129").Value = 3
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 'H'
$ws.Range("K129").Value = 1.363
$ws.Range("L129").Value = 4.75
$ws.Range("M129").Value = 6.5
$ws.Range("N129").Value = 1.285
$ws.Range("O129").Value = 6.5
$ws.Range("P129").Value = 8
$ws.Range("Q129").Value = -1.75
$ws.Range("R129").Value = 1.9
$ws.Range("S129").Value = 1.9
$ws.Range("T129").Value = 3.25
$ws.Range("U129").Value = 1.85
$ws.Range("V129").Value = 1.95
$ws.Range("W129").Value = 0.2849999999999999
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.8999999999999999
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = -0.5
$ws.Range("AC129").Value = 0.475

# Row 148
$ws.Range("B148").Value = 7532421
$ws.Range("F148").Value = 'Guabira'
$ws.Range("G148").Value = 'Independiente Petrolero'
$ws.Range("H148").Value = 2
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = 'H'
$ws.Range("K148").Value = 1.4
$ws.Range("L148").Value = 4.5
$ws.Range("M148").Value = 7.5
$ws.Range("N148").Value = 1.333
$ws.Range("O148").Value = 5.5
$ws.Range("P148").Value = 9.5
$ws.Range("Q148").Value = -1.5
$ws.Range("R148").Value = 1.85
$ws.Range("S148").Value = 1.95
$ws.Range("T148").Value = 3
$ws.Range("U148").Value = 1.825
$ws.Range("V148").Value = 1.975
$ws.Range("W148").Value = 0.333
$ws.Range("X148").Value = -1
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = 0.8500000000000001
$ws.Range("AA148").Value = -1
$ws.Range("AB148").Value = -1
$ws.Range("AC148").Value = 0.9750000000000001

# Row 149
$ws.Range("B149").Value = 7532419
$ws.Range("F149").Value = 'Oriente Petrolero'
$ws.Range("G149").Value = 'Jorge Wilstermann'
$ws.Range("H149").Value = 3
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 'H'
$ws.Range("K149").Value = 2.2
$ws.Range("L149").Value = 2.5
$ws.Range("M149").Value = 4.5
$ws.Range("N149").Value = 2.375
$ws.Range("O149").Value = 2.45
$ws.Range("P149").Value = 4.5
$ws.Range("Q149").Value = -0.25
$ws.Range("R149").Value = 1.9
$ws.Range("S149").Value = 1.9
$ws.Range("T149").Value = 2
$ws.Range("U149").Value = 1.95
$ws.Range("V149").Value = 1.85
$ws.Range("W149").Value = 1.375
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 0.8999999999999999
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = 0.95
$ws.Range("AC149").Value = -1

# Row 150
$ws.Range("B150").Value = 7532420
$ws.Range("F150").Value = 'Club Aurora'
$ws.Range("G150").Value = 'Vaca Diez'
$ws.Range("H150").Value = 3
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 'H'
$ws.Range("K150").Value = 1.333
$ws.Range("L150").Value = 5
$ws.Range("M150").Value = 8
$ws.Range("N150").Value = 1.3
$ws.Range("O150").Value = 6.5
$ws.Range("P150").Value = 7
$ws.Range("Q150").Value = -1.5
$ws.Range("R150").Value = 1.8
$ws.Range("S150").Value = 2
$ws.Range("T150").Value = 3.25
$ws.Range("U150").Value = 1.95
$ws.Range("V150").Value = 1.85
$ws.Range("W150").Value = 0.3
$ws.Range("X150").Value = -1
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = 0.8
$ws.Range("AA150").Value = -1
$ws.Range("AB150").Value = -0.5
$ws.Range("AC150").Value = 0.425
